$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells - copy formatting from an existing header cell (H1) so the
# new cells end up sharing the same style definition as the rest of row 1.
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for columns I (I0) and J (IF), rows 2-30
$data = @(
    @(8, 9),
    @(8, 8),
    @(5, 5),
    @(5, 6),
    @(9, 9),
    @(7, 7),
    @(6, 6),
    @(6, 6),
    @(4, 4),
    @(7, 7),
    @(8, 8),
    @(7, 7),
    @(8, 8),
    @(7, 7),
    @(6, 6),
    @(6, 6),
    @(6, 6),
    @(7, 7),
    @(6, 6),
    @(5, 6),
    @(6, 7),
    @(6, 6),
    @(8, 8),
    @(5, 5),
    @(6, 6),
    @(8, 9),
    @(9, 9),
    @(6, 6),
    @(5, 5)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
